# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) on each data sheet, drops the
# bold/bordered style from the former "header-ish" A-column label cells
# (rows below the new header row), fixes a handful of accented labels,
# removes the now-unused "Teto" row on the Emissoes sheet, and refreshes
# the Custo Total sheet with a header row + updated figures.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# ---------------------------------------------------------------------
$dataSheets = @(1, 2, 3, 4)

foreach ($idx in $dataSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # New header cell in A1, formatted like the other header cells (B1:E1).
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value2 = "Fonte/Tecnologia"

    # Strip the bold/border style from the row-label cells (A2:A12) -
    # they keep plain/default formatting now that A1 carries the header.
    $ws.Range("A2:A12").Style = "Normal"

    # Fix accented labels.
    $ws.Range("A3").Value2 = "Gás Natural"
    $ws.Range("A4").Value2 = "Carvão"
    $ws.Range("A6").Value2 = "Óleos Comb"
    $ws.Range("A8").Value2 = "Eólica"
    $ws.Range("A11").Value2 = "Pot. Compl."
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value2 = "Período"

$ws5.Range("A2:A3").Style = "Normal"
$ws5.Range("A2").Value2 = "P.Médio"
$ws5.Range("A3").Value2 = "P.Crítico"

# Remove the "Teto" row entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value2 = "Tipo Expansão"

# "2015" must stay text (matches the other sheets' text-typed year headers
# like B1 on sheet 1) rather than becoming a real number, which is what a
# plain `Value2 = "2015"` assignment would auto-coerce it into. Pasting the
# *value* of an existing text-typed "2015" cell keeps the literal text type
# while leaving B1's own style (already bold/bordered) untouched.
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4163)

$ws6.Range("A2:A3").Style = "Normal"
$ws6.Range("A2").Value2 = "Expansão Centralizada"
$ws6.Range("B2").Value2 = 595
$ws6.Range("A3").Value2 = "Expansão por GD"
$ws6.Range("B3").Value2 = 99
